# Bugfixed QoQ Visualizations and a typo in the evaluation objects
# Removes the first 16 data rows (rows 2-17) from the sheet, shifting all
# subsequent rows up so that the data starts from the same series but
# 16 rows later (new row 2 == old row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B17").EntireRow.Delete() | Out-Null
